$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H26").Value = 0.6514799999999999
$ws.Range("I26").Value = 0.03403
$ws.Range("H27").Value = 0.08021
$ws.Range("I27").Value = 0.05101
$ws.Range("H28").Value = 0.65222
$ws.Range("I28").Value = 0.03476
$ws.Range("H29").Value = 0.11531
$ws.Range("I29").Value = 0.06062
$ws.Range("H30").Value = 0.65166
$ws.Range("I30").Value = 0.03342
$ws.Range("H31").Value = 0.08214
$ws.Range("I31").Value = 0.05028
$ws.Range("H32").Value = 0.65237
$ws.Range("I32").Value = 0.03555
$ws.Range("H33").Value = 0.12309
$ws.Range("I33").Value = 0.06365999999999999
$ws.Range("H34").Value = 0.65679
$ws.Range("I34").Value = 0.02357
$ws.Range("H35").Value = 0.04355
$ws.Range("I35").Value = 0.03496
$ws.Range("H36").Value = 0.6577499999999999
$ws.Range("I36").Value = 0.0234
$ws.Range("H37").Value = 0.04841
$ws.Range("I37").Value = 0.03818
$ws.Range("H38").Value = 0.65691
$ws.Range("I38").Value = 0.02339
$ws.Range("H39").Value = 0.04355
$ws.Range("I39").Value = 0.03496
$ws.Range("H40").Value = 0.65782
$ws.Range("I40").Value = 0.02357
$ws.Range("H41").Value = 0.04938
$ws.Range("I41").Value = 0.0391
$ws.Range("H66").Value = 0.61674
$ws.Range("I66").Value = 0.03628
$ws.Range("H67").Value = 0.05807
$ws.Range("I67").Value = 0.04235
$ws.Range("H68").Value = 0.62239
$ws.Range("I68").Value = 0.03595
$ws.Range("H69").Value = 0.07649
$ws.Range("I69").Value = 0.05531
$ws.Range("H70").Value = 0.61826
$ws.Range("I70").Value = 0.03455
$ws.Range("H71").Value = 0.05805
$ws.Range("I71").Value = 0.04123
$ws.Range("H72").Value = 0.62485
$ws.Range("I72").Value = 0.03711
$ws.Range("H73").Value = 0.09296
$ws.Range("I73").Value = 0.05759
$ws.Range("H74").Value = 0.6317700000000001
$ws.Range("I74").Value = 0.03569
$ws.Range("H75").Value = 0.03289
$ws.Range("I75").Value = 0.03146
$ws.Range("H76").Value = 0.6359
$ws.Range("I76").Value = 0.03562
$ws.Range("H77").Value = 0.04358
$ws.Range("I77").Value = 0.03837
$ws.Range("H78").Value = 0.63236
$ws.Range("I78").Value = 0.03578
$ws.Range("H79").Value = 0.03289
$ws.Range("I79").Value = 0.03146
$ws.Range("H80").Value = 0.63596
$ws.Range("I80").Value = 0.03611
$ws.Range("H81").Value = 0.04455
$ws.Range("I81").Value = 0.03941
$ws.Range("H106").Value = 0.65344
$ws.Range("I106").Value = 0.02625
$ws.Range("H107").Value = 0.09868
$ws.Range("I107").Value = 0.05485
$ws.Range("H108").Value = 0.6539700000000001
$ws.Range("I108").Value = 0.02843
$ws.Range("H109").Value = 0.12483
$ws.Range("I109").Value = 0.05707
$ws.Range("H110").Value = 0.65422
$ws.Range("I110").Value = 0.02601
$ws.Range("H111").Value = 0.09868
$ws.Range("I111").Value = 0.05441
$ws.Range("H112").Value = 0.65572
$ws.Range("I112").Value = 0.02778
$ws.Range("H113").Value = 0.13359
$ws.Range("I113").Value = 0.05731
$ws.Range("H114").Value = 0.6533
$ws.Range("I114").Value = 0.02206
$ws.Range("H115").Value = 0.03877
$ws.Range("I115").Value = 0.03624
$ws.Range("H116").Value = 0.65452
$ws.Range("I116").Value = 0.02151
$ws.Range("H117").Value = 0.0426
$ws.Range("I117").Value = 0.03765
$ws.Range("H118").Value = 0.65324
$ws.Range("I118").Value = 0.02212
$ws.Range("H119").Value = 0.03877
$ws.Range("I119").Value = 0.03624
$ws.Range("H120").Value = 0.6541400000000001
$ws.Range("I120").Value = 0.02181
$ws.Range("H121").Value = 0.0426
$ws.Range("I121").Value = 0.03765
$ws.Range("H146").Value = 0.63334
$ws.Range("I146").Value = 0.02814
$ws.Range("H147").Value = 0.05612
$ws.Range("I147").Value = 0.0381
$ws.Range("H148").Value = 0.63464
$ws.Range("I148").Value = 0.0308
$ws.Range("H149").Value = 0.09085
$ws.Range("I149").Value = 0.04875
$ws.Range("H150").Value = 0.6339900000000001
$ws.Range("I150").Value = 0.02775
$ws.Range("H151").Value = 0.05612
$ws.Range("I151").Value = 0.0381
$ws.Range("H152").Value = 0.63898
$ws.Range("I152").Value = 0.03177
$ws.Range("H153").Value = 0.11501
$ws.Range("I153").Value = 0.04646
$ws.Range("H154").Value = 0.64498
$ws.Range("I154").Value = 0.02504
$ws.Range("H155").Value = 0.03199
$ws.Range("I155").Value = 0.03128
$ws.Range("H156").Value = 0.64727
$ws.Range("I156").Value = 0.02511
$ws.Range("H157").Value = 0.04643
$ws.Range("I157").Value = 0.03744
$ws.Range("H158").Value = 0.64479
$ws.Range("I158").Value = 0.02512
$ws.Range("H159").Value = 0.03199
$ws.Range("I159").Value = 0.03128
$ws.Range("H160").Value = 0.6476499999999999
$ws.Range("I160").Value = 0.0254
$ws.Range("H161").Value = 0.04839
$ws.Range("I161").Value = 0.03814
